# Auto-generated edit script for workbook update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 102-237 (weekly data shifted down by one week; new week inserted at row 102) ---
$ws.Cells.Item(102, 4).Value = 44579
$ws.Cells.Item(102, 10).Value = 800
$ws.Cells.Item(102, 11).Value = 1000
$ws.Cells.Item(102, 12).Value = 1000
$ws.Cells.Item(102, 13).Value = 1000
$ws.Cells.Item(102, 16).Value = 200
$ws.Cells.Item(103, 4).Value = 44284
$ws.Cells.Item(104, 4).Value = 44175
$ws.Cells.Item(104, 11).Value = 850
$ws.Cells.Item(104, 12).Value = 850
$ws.Cells.Item(104, 13).Value = 850
$ws.Cells.Item(104, 16).Value = 170
$ws.Cells.Item(105, 4).Value = 44242
$ws.Cells.Item(105, 10).Value = 500
$ws.Cells.Item(106, 4).Value = 44419
$ws.Cells.Item(106, 10).Value = 250
$ws.Cells.Item(107, 4).Value = 44424
$ws.Cells.Item(108, 4).Value = 44202
$ws.Cells.Item(108, 10).Value = 500
$ws.Cells.Item(108, 11).Value = 1000
$ws.Cells.Item(108, 13).Value = 1000
$ws.Cells.Item(108, 16).Value = 200
$ws.Cells.Item(109, 4).Value = 44271
$ws.Cells.Item(109, 10).Value = 1000
$ws.Cells.Item(109, 11).Value = 800
$ws.Cells.Item(109, 13).Value = 900
$ws.Cells.Item(109, 16).Value = 180
$ws.Cells.Item(110, 4).Value = 44364
$ws.Cells.Item(111, 4).Value = 44259
$ws.Cells.Item(111, 10).Value = 500
$ws.Cells.Item(111, 11).Value = 1000
$ws.Cells.Item(111, 13).Value = 1000
$ws.Cells.Item(111, 16).Value = 200
$ws.Cells.Item(112, 4).Value = 44309
$ws.Cells.Item(112, 10).Value = 1200
$ws.Cells.Item(112, 11).Value = 900
$ws.Cells.Item(112, 13).Value = 950
$ws.Cells.Item(112, 16).Value = 190
$ws.Cells.Item(113, 4).Value = 44508
$ws.Cells.Item(113, 10).Value = 500
$ws.Cells.Item(113, 11).Value = 1000
$ws.Cells.Item(113, 13).Value = 1000
$ws.Cells.Item(113, 16).Value = 200
$ws.Cells.Item(114, 4).Value = 44540
$ws.Cells.Item(114, 10).Value = 1200
$ws.Cells.Item(114, 11).Value = 900
$ws.Cells.Item(114, 13).Value = 950
$ws.Cells.Item(114, 16).Value = 190
$ws.Cells.Item(115, 4).Value = 44391
$ws.Cells.Item(116, 4).Value = 44417
$ws.Cells.Item(116, 11).Value = 1000
$ws.Cells.Item(116, 13).Value = 1000
$ws.Cells.Item(116, 16).Value = 200
$ws.Cells.Item(117, 4).Value = 44273
$ws.Cells.Item(117, 11).Value = 850
$ws.Cells.Item(117, 13).Value = 925
$ws.Cells.Item(117, 16).Value = 185
$ws.Cells.Item(118, 4).Value = 44529
$ws.Cells.Item(118, 11).Value = 1000
$ws.Cells.Item(118, 13).Value = 1000
$ws.Cells.Item(118, 16).Value = 200
$ws.Cells.Item(119, 4).Value = 44214
$ws.Cells.Item(119, 11).Value = 900
$ws.Cells.Item(119, 13).Value = 950
$ws.Cells.Item(119, 16).Value = 190
$ws.Cells.Item(120, 4).Value = 44567
$ws.Cells.Item(120, 10).Value = 500
$ws.Cells.Item(121, 4).Value = 44421
$ws.Cells.Item(121, 10).Value = 1200
$ws.Cells.Item(121, 11).Value = 1000
$ws.Cells.Item(121, 13).Value = 1000
$ws.Cells.Item(121, 16).Value = 200
$ws.Cells.Item(122, 4).Value = 44546
$ws.Cells.Item(122, 10).Value = 500
$ws.Cells.Item(122, 11).Value = 800
$ws.Cells.Item(122, 12).Value = 1000
$ws.Cells.Item(122, 13).Value = 900
$ws.Cells.Item(122, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(122, 15).Value = 'Región del Maule'
$ws.Cells.Item(122, 16).Value = 180
$ws.Cells.Item(122, 17).Value = 5
$ws.Cells.Item(123, 10).Value = 110
$ws.Cells.Item(123, 11).Value = 8000
$ws.Cells.Item(123, 12).Value = 8000
$ws.Cells.Item(123, 13).Value = 8000
$ws.Cells.Item(123, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(123, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(123, 16).Value = 533
$ws.Cells.Item(123, 17).Value = 15
$ws.Cells.Item(124, 4).Value = 44400
$ws.Cells.Item(124, 10).Value = 1200
$ws.Cells.Item(125, 4).Value = 44350
$ws.Cells.Item(125, 11).Value = 1000
$ws.Cells.Item(125, 12).Value = 1000
$ws.Cells.Item(125, 13).Value = 1000
$ws.Cells.Item(125, 16).Value = 200
$ws.Cells.Item(126, 4).Value = 44172
$ws.Cells.Item(126, 10).Value = 500
$ws.Cells.Item(126, 11).Value = 800
$ws.Cells.Item(126, 12).Value = 800
$ws.Cells.Item(126, 13).Value = 800
$ws.Cells.Item(126, 16).Value = 160
$ws.Cells.Item(127, 4).Value = 44475
$ws.Cells.Item(127, 10).Value = 250
$ws.Cells.Item(127, 11).Value = 1000
$ws.Cells.Item(127, 12).Value = 1000
$ws.Cells.Item(127, 13).Value = 1000
$ws.Cells.Item(127, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(127, 15).Value = 'Región del Maule'
$ws.Cells.Item(127, 16).Value = 200
$ws.Cells.Item(127, 17).Value = 5
$ws.Cells.Item(128, 10).Value = 110
$ws.Cells.Item(128, 11).Value = 8000
$ws.Cells.Item(128, 12).Value = 8000
$ws.Cells.Item(128, 13).Value = 8000
$ws.Cells.Item(128, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(128, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(128, 16).Value = 533
$ws.Cells.Item(128, 17).Value = 15
$ws.Cells.Item(129, 4).Value = 44201
$ws.Cells.Item(129, 10).Value = 1200
$ws.Cells.Item(129, 11).Value = 900
$ws.Cells.Item(129, 12).Value = 900
$ws.Cells.Item(130, 4).Value = 44267
$ws.Cells.Item(130, 10).Value = 1000
$ws.Cells.Item(130, 11).Value = 800
$ws.Cells.Item(130, 12).Value = 1000
$ws.Cells.Item(130, 13).Value = 900
$ws.Cells.Item(130, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(130, 15).Value = 'Región del Maule'
$ws.Cells.Item(130, 16).Value = 180
$ws.Cells.Item(130, 17).Value = 5
$ws.Cells.Item(131, 10).Value = 150
$ws.Cells.Item(131, 11).Value = 10000
$ws.Cells.Item(131, 12).Value = 10000
$ws.Cells.Item(131, 13).Value = 10000
$ws.Cells.Item(131, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(131, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(131, 16).Value = 667
$ws.Cells.Item(131, 17).Value = 15
$ws.Cells.Item(132, 4).Value = 44215
$ws.Cells.Item(132, 10).Value = 1200
$ws.Cells.Item(133, 4).Value = 44539
$ws.Cells.Item(133, 10).Value = 500
$ws.Cells.Item(133, 11).Value = 900
$ws.Cells.Item(133, 13).Value = 950
$ws.Cells.Item(133, 16).Value = 190
$ws.Cells.Item(134, 4).Value = 44319
$ws.Cells.Item(134, 10).Value = 250
$ws.Cells.Item(134, 11).Value = 1000
$ws.Cells.Item(134, 13).Value = 1000
$ws.Cells.Item(134, 16).Value = 200
$ws.Cells.Item(135, 4).Value = 44523
$ws.Cells.Item(135, 11).Value = 900
$ws.Cells.Item(135, 12).Value = 1000
$ws.Cells.Item(135, 13).Value = 950
$ws.Cells.Item(135, 16).Value = 190
$ws.Cells.Item(136, 4).Value = 44495
$ws.Cells.Item(136, 10).Value = 1200
$ws.Cells.Item(136, 11).Value = 1000
$ws.Cells.Item(136, 12).Value = 1200
$ws.Cells.Item(136, 13).Value = 1100
$ws.Cells.Item(136, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(136, 15).Value = 'Región del Maule'
$ws.Cells.Item(136, 16).Value = 220
$ws.Cells.Item(136, 17).Value = 5
$ws.Cells.Item(137, 10).Value = 120
$ws.Cells.Item(137, 11).Value = 7000
$ws.Cells.Item(137, 12).Value = 8000
$ws.Cells.Item(137, 13).Value = 7500
$ws.Cells.Item(137, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(137, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(137, 16).Value = 500
$ws.Cells.Item(137, 17).Value = 15
$ws.Cells.Item(138, 4).Value = 44292
$ws.Cells.Item(138, 10).Value = 1000
$ws.Cells.Item(138, 11).Value = 850
$ws.Cells.Item(138, 12).Value = 1000
$ws.Cells.Item(138, 13).Value = 925
$ws.Cells.Item(138, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(138, 15).Value = 'Región del Maule'
$ws.Cells.Item(138, 16).Value = 185
$ws.Cells.Item(138, 17).Value = 5
$ws.Cells.Item(139, 10).Value = 110
$ws.Cells.Item(139, 11).Value = 10000
$ws.Cells.Item(139, 12).Value = 10000
$ws.Cells.Item(139, 13).Value = 10000
$ws.Cells.Item(139, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(139, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(139, 16).Value = 667
$ws.Cells.Item(139, 17).Value = 15
$ws.Cells.Item(140, 4).Value = 44204
$ws.Cells.Item(140, 10).Value = 1000
$ws.Cells.Item(140, 11).Value = 900
$ws.Cells.Item(140, 12).Value = 900
$ws.Cells.Item(140, 13).Value = 900
$ws.Cells.Item(140, 16).Value = 180
$ws.Cells.Item(141, 4).Value = 44574
$ws.Cells.Item(141, 11).Value = 1000
$ws.Cells.Item(141, 12).Value = 1000
$ws.Cells.Item(141, 13).Value = 1000
$ws.Cells.Item(141, 16).Value = 200
$ws.Cells.Item(142, 4).Value = 44168
$ws.Cells.Item(142, 11).Value = 850
$ws.Cells.Item(142, 12).Value = 850
$ws.Cells.Item(142, 13).Value = 850
$ws.Cells.Item(142, 16).Value = 170
$ws.Cells.Item(143, 4).Value = 44179
$ws.Cells.Item(143, 10).Value = 500
$ws.Cells.Item(144, 4).Value = 44299
$ws.Cells.Item(144, 10).Value = 1250
$ws.Cells.Item(144, 11).Value = 1000
$ws.Cells.Item(144, 13).Value = 1000
$ws.Cells.Item(144, 16).Value = 200
$ws.Cells.Item(145, 4).Value = 44161
$ws.Cells.Item(145, 10).Value = 500
$ws.Cells.Item(145, 11).Value = 950
$ws.Cells.Item(145, 12).Value = 1000
$ws.Cells.Item(145, 13).Value = 975
$ws.Cells.Item(145, 16).Value = 195
$ws.Cells.Item(146, 4).Value = 44477
$ws.Cells.Item(146, 10).Value = 1200
$ws.Cells.Item(146, 12).Value = 1200
$ws.Cells.Item(146, 13).Value = 1100
$ws.Cells.Item(146, 16).Value = 220
$ws.Cells.Item(147, 4).Value = 44438
$ws.Cells.Item(147, 10).Value = 500
$ws.Cells.Item(147, 11).Value = 1000
$ws.Cells.Item(147, 13).Value = 1000
$ws.Cells.Item(147, 16).Value = 200
$ws.Cells.Item(148, 4).Value = 44498
$ws.Cells.Item(148, 10).Value = 1400
$ws.Cells.Item(148, 11).Value = 900
$ws.Cells.Item(148, 13).Value = 950
$ws.Cells.Item(148, 16).Value = 190
$ws.Cells.Item(149, 4).Value = 44413
$ws.Cells.Item(149, 10).Value = 500
$ws.Cells.Item(149, 11).Value = 1000
$ws.Cells.Item(149, 13).Value = 1000
$ws.Cells.Item(149, 16).Value = 200
$ws.Cells.Item(150, 4).Value = 44544
$ws.Cells.Item(150, 10).Value = 1200
$ws.Cells.Item(150, 11).Value = 900
$ws.Cells.Item(150, 12).Value = 1000
$ws.Cells.Item(150, 13).Value = 950
$ws.Cells.Item(150, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(150, 15).Value = 'Región del Maule'
$ws.Cells.Item(150, 16).Value = 190
$ws.Cells.Item(150, 17).Value = 5
$ws.Cells.Item(151, 4).Value = 44160
$ws.Cells.Item(151, 10).Value = 20
$ws.Cells.Item(151, 11).Value = 8000
$ws.Cells.Item(151, 12).Value = 8000
$ws.Cells.Item(151, 13).Value = 8000
$ws.Cells.Item(151, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(151, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(151, 16).Value = 533
$ws.Cells.Item(151, 17).Value = 15
$ws.Cells.Item(152, 4).Value = 44357
$ws.Cells.Item(152, 10).Value = 400
$ws.Cells.Item(152, 11).Value = 1000
$ws.Cells.Item(152, 12).Value = 1000
$ws.Cells.Item(152, 13).Value = 1000
$ws.Cells.Item(152, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(152, 15).Value = 'Región del Maule'
$ws.Cells.Item(152, 16).Value = 200
$ws.Cells.Item(152, 17).Value = 5
$ws.Cells.Item(153, 10).Value = 120
$ws.Cells.Item(153, 11).Value = 9000
$ws.Cells.Item(153, 12).Value = 9000
$ws.Cells.Item(153, 13).Value = 9000
$ws.Cells.Item(153, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(153, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(153, 16).Value = 600
$ws.Cells.Item(153, 17).Value = 15
$ws.Cells.Item(154, 4).Value = 44218
$ws.Cells.Item(154, 10).Value = 1000
$ws.Cells.Item(154, 11).Value = 800
$ws.Cells.Item(154, 12).Value = 900
$ws.Cells.Item(154, 13).Value = 850
$ws.Cells.Item(154, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(154, 15).Value = 'Región del Maule'
$ws.Cells.Item(154, 16).Value = 170
$ws.Cells.Item(154, 17).Value = 5
$ws.Cells.Item(155, 10).Value = 110
$ws.Cells.Item(155, 11).Value = 10000
$ws.Cells.Item(155, 12).Value = 10000
$ws.Cells.Item(155, 13).Value = 10000
$ws.Cells.Item(155, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(155, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(155, 16).Value = 667
$ws.Cells.Item(155, 17).Value = 15
$ws.Cells.Item(156, 4).Value = 44320
$ws.Cells.Item(157, 4).Value = 44306
$ws.Cells.Item(158, 4).Value = 44295
$ws.Cells.Item(158, 10).Value = 1200
$ws.Cells.Item(158, 11).Value = 1000
$ws.Cells.Item(158, 13).Value = 1000
$ws.Cells.Item(158, 16).Value = 200
$ws.Cells.Item(159, 4).Value = 44210
$ws.Cells.Item(159, 10).Value = 750
$ws.Cells.Item(159, 11).Value = 900
$ws.Cells.Item(159, 13).Value = 967
$ws.Cells.Item(159, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(159, 16).Value = 193
$ws.Cells.Item(159, 17).Value = 5
$ws.Cells.Item(160, 4).Value = 44407
$ws.Cells.Item(160, 10).Value = 1200
$ws.Cells.Item(160, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(160, 16).Value = 67
$ws.Cells.Item(160, 17).Value = 15
$ws.Cells.Item(161, 4).Value = 44343
$ws.Cells.Item(161, 10).Value = 500
$ws.Cells.Item(162, 4).Value = 44230
$ws.Cells.Item(162, 10).Value = 250
$ws.Cells.Item(162, 11).Value = 1000
$ws.Cells.Item(162, 12).Value = 1000
$ws.Cells.Item(162, 13).Value = 1000
$ws.Cells.Item(162, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(162, 15).Value = 'Región del Maule'
$ws.Cells.Item(162, 16).Value = 200
$ws.Cells.Item(162, 17).Value = 5
$ws.Cells.Item(163, 10).Value = 110
$ws.Cells.Item(163, 11).Value = 8000
$ws.Cells.Item(163, 12).Value = 8000
$ws.Cells.Item(163, 13).Value = 8000
$ws.Cells.Item(163, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(163, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(163, 16).Value = 533
$ws.Cells.Item(163, 17).Value = 15
$ws.Cells.Item(164, 4).Value = 44316
$ws.Cells.Item(164, 10).Value = 1200
$ws.Cells.Item(165, 4).Value = 44265
$ws.Cells.Item(165, 10).Value = 250
$ws.Cells.Item(166, 4).Value = 44329
$ws.Cells.Item(167, 4).Value = 44460
$ws.Cells.Item(167, 10).Value = 500
$ws.Cells.Item(167, 11).Value = 1000
$ws.Cells.Item(167, 13).Value = 1000
$ws.Cells.Item(167, 16).Value = 200
$ws.Cells.Item(168, 4).Value = 44526
$ws.Cells.Item(168, 10).Value = 1200
$ws.Cells.Item(168, 11).Value = 900
$ws.Cells.Item(168, 12).Value = 1000
$ws.Cells.Item(168, 13).Value = 950
$ws.Cells.Item(168, 16).Value = 190
$ws.Cells.Item(169, 4).Value = 44467
$ws.Cells.Item(169, 10).Value = 1100
$ws.Cells.Item(169, 11).Value = 1200
$ws.Cells.Item(169, 12).Value = 1200
$ws.Cells.Item(169, 13).Value = 1200
$ws.Cells.Item(169, 16).Value = 240
$ws.Cells.Item(170, 4).Value = 44411
$ws.Cells.Item(170, 10).Value = 1200
$ws.Cells.Item(170, 11).Value = 1000
$ws.Cells.Item(170, 12).Value = 1000
$ws.Cells.Item(170, 13).Value = 1000
$ws.Cells.Item(170, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(170, 15).Value = 'Región del Maule'
$ws.Cells.Item(170, 16).Value = 200
$ws.Cells.Item(170, 17).Value = 5
$ws.Cells.Item(171, 10).Value = 120
$ws.Cells.Item(171, 11).Value = 10000
$ws.Cells.Item(171, 12).Value = 10000
$ws.Cells.Item(171, 13).Value = 10000
$ws.Cells.Item(171, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(171, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(171, 16).Value = 667
$ws.Cells.Item(171, 17).Value = 15
$ws.Cells.Item(172, 4).Value = 44211
$ws.Cells.Item(172, 11).Value = 900
$ws.Cells.Item(172, 13).Value = 950
$ws.Cells.Item(172, 16).Value = 190
$ws.Cells.Item(173, 4).Value = 44313
$ws.Cells.Item(173, 10).Value = 1000
$ws.Cells.Item(173, 11).Value = 1000
$ws.Cells.Item(173, 12).Value = 1000
$ws.Cells.Item(173, 13).Value = 1000
$ws.Cells.Item(173, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(173, 15).Value = 'Región del Maule'
$ws.Cells.Item(173, 16).Value = 200
$ws.Cells.Item(173, 17).Value = 5
$ws.Cells.Item(174, 10).Value = 120
$ws.Cells.Item(174, 11).Value = 9000
$ws.Cells.Item(174, 12).Value = 9000
$ws.Cells.Item(174, 13).Value = 9000
$ws.Cells.Item(174, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(174, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(174, 16).Value = 600
$ws.Cells.Item(174, 17).Value = 15
$ws.Cells.Item(175, 4).Value = 44334
$ws.Cells.Item(175, 10).Value = 1200
$ws.Cells.Item(175, 11).Value = 1000
$ws.Cells.Item(175, 12).Value = 1000
$ws.Cells.Item(175, 13).Value = 1000
$ws.Cells.Item(175, 16).Value = 200
$ws.Cells.Item(176, 4).Value = 44209
$ws.Cells.Item(176, 10).Value = 500
$ws.Cells.Item(176, 11).Value = 900
$ws.Cells.Item(176, 12).Value = 900
$ws.Cells.Item(176, 13).Value = 900
$ws.Cells.Item(176, 16).Value = 180
$ws.Cells.Item(177, 4).Value = 44389
$ws.Cells.Item(177, 10).Value = 250
$ws.Cells.Item(177, 11).Value = 1000
$ws.Cells.Item(177, 12).Value = 1000
$ws.Cells.Item(177, 13).Value = 1000
$ws.Cells.Item(177, 16).Value = 200
$ws.Cells.Item(178, 4).Value = 44221
$ws.Cells.Item(178, 11).Value = 900
$ws.Cells.Item(178, 12).Value = 900
$ws.Cells.Item(178, 13).Value = 900
$ws.Cells.Item(178, 16).Value = 180
$ws.Cells.Item(179, 4).Value = 44280
$ws.Cells.Item(179, 10).Value = 500
$ws.Cells.Item(179, 11).Value = 850
$ws.Cells.Item(179, 12).Value = 1000
$ws.Cells.Item(179, 13).Value = 925
$ws.Cells.Item(179, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(179, 15).Value = 'Región del Maule'
$ws.Cells.Item(179, 16).Value = 185
$ws.Cells.Item(179, 17).Value = 5
$ws.Cells.Item(180, 10).Value = 110
$ws.Cells.Item(180, 11).Value = 10000
$ws.Cells.Item(180, 12).Value = 10000
$ws.Cells.Item(180, 13).Value = 10000
$ws.Cells.Item(180, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(180, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(180, 16).Value = 667
$ws.Cells.Item(180, 17).Value = 15
$ws.Cells.Item(181, 4).Value = 44330
$ws.Cells.Item(181, 10).Value = 1200
$ws.Cells.Item(182, 4).Value = 44483
$ws.Cells.Item(183, 4).Value = 44448
$ws.Cells.Item(183, 10).Value = 500
$ws.Cells.Item(183, 11).Value = 1000
$ws.Cells.Item(183, 12).Value = 1000
$ws.Cells.Item(183, 13).Value = 1000
$ws.Cells.Item(183, 16).Value = 200
$ws.Cells.Item(184, 4).Value = 44463
$ws.Cells.Item(184, 11).Value = 1200
$ws.Cells.Item(184, 12).Value = 1200
$ws.Cells.Item(184, 13).Value = 1200
$ws.Cells.Item(184, 16).Value = 240
$ws.Cells.Item(185, 4).Value = 44239
$ws.Cells.Item(185, 10).Value = 1200
$ws.Cells.Item(185, 12).Value = 1000
$ws.Cells.Item(185, 13).Value = 1000
$ws.Cells.Item(185, 16).Value = 200
$ws.Cells.Item(186, 4).Value = 44476
$ws.Cells.Item(186, 10).Value = 500
$ws.Cells.Item(186, 11).Value = 1000
$ws.Cells.Item(186, 12).Value = 1200
$ws.Cells.Item(186, 13).Value = 1100
$ws.Cells.Item(186, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(186, 15).Value = 'Región del Maule'
$ws.Cells.Item(186, 16).Value = 220
$ws.Cells.Item(186, 17).Value = 5
$ws.Cells.Item(187, 10).Value = 120
$ws.Cells.Item(187, 11).Value = 8000
$ws.Cells.Item(187, 12).Value = 8000
$ws.Cells.Item(187, 13).Value = 8000
$ws.Cells.Item(187, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(187, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(187, 16).Value = 533
$ws.Cells.Item(187, 17).Value = 15
$ws.Cells.Item(188, 4).Value = 44169
$ws.Cells.Item(188, 10).Value = 1200
$ws.Cells.Item(188, 11).Value = 850
$ws.Cells.Item(188, 12).Value = 1000
$ws.Cells.Item(188, 13).Value = 925
$ws.Cells.Item(188, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(188, 15).Value = 'Región del Maule'
$ws.Cells.Item(188, 16).Value = 185
$ws.Cells.Item(188, 17).Value = 5
$ws.Cells.Item(189, 10).Value = 120
$ws.Cells.Item(189, 11).Value = 10000
$ws.Cells.Item(189, 12).Value = 10000
$ws.Cells.Item(189, 13).Value = 10000
$ws.Cells.Item(189, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(189, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(189, 16).Value = 667
$ws.Cells.Item(189, 17).Value = 15
$ws.Cells.Item(190, 4).Value = 44208
$ws.Cells.Item(190, 12).Value = 900
$ws.Cells.Item(190, 13).Value = 900
$ws.Cells.Item(190, 16).Value = 180
$ws.Cells.Item(191, 4).Value = 44250
$ws.Cells.Item(191, 10).Value = 1200
$ws.Cells.Item(191, 13).Value = 950
$ws.Cells.Item(191, 16).Value = 190
$ws.Cells.Item(192, 4).Value = 44515
$ws.Cells.Item(192, 10).Value = 750
$ws.Cells.Item(192, 11).Value = 900
$ws.Cells.Item(192, 12).Value = 1000
$ws.Cells.Item(192, 13).Value = 933
$ws.Cells.Item(192, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(192, 15).Value = 'Región del Maule'
$ws.Cells.Item(192, 16).Value = 187
$ws.Cells.Item(192, 17).Value = 5
$ws.Cells.Item(193, 10).Value = 120
$ws.Cells.Item(193, 11).Value = 9000
$ws.Cells.Item(193, 12).Value = 9000
$ws.Cells.Item(193, 13).Value = 9000
$ws.Cells.Item(193, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(193, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(193, 16).Value = 600
$ws.Cells.Item(193, 17).Value = 15
$ws.Cells.Item(194, 4).Value = 44379
$ws.Cells.Item(194, 10).Value = 1200
$ws.Cells.Item(194, 11).Value = 1000
$ws.Cells.Item(194, 12).Value = 1000
$ws.Cells.Item(194, 13).Value = 1000
$ws.Cells.Item(194, 16).Value = 200
$ws.Cells.Item(195, 4).Value = 44216
$ws.Cells.Item(195, 10).Value = 250
$ws.Cells.Item(195, 11).Value = 900
$ws.Cells.Item(195, 12).Value = 900
$ws.Cells.Item(195, 13).Value = 900
$ws.Cells.Item(195, 16).Value = 180
$ws.Cells.Item(196, 4).Value = 44509
$ws.Cells.Item(196, 10).Value = 1200
$ws.Cells.Item(196, 12).Value = 1200
$ws.Cells.Item(196, 13).Value = 1100
$ws.Cells.Item(196, 15).Value = 'Región del Maule'
$ws.Cells.Item(196, 16).Value = 220
$ws.Cells.Item(197, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(198, 4).Value = 44229
$ws.Cells.Item(198, 10).Value = 500
$ws.Cells.Item(198, 12).Value = 1000
$ws.Cells.Item(198, 13).Value = 1000
$ws.Cells.Item(198, 16).Value = 200
$ws.Cells.Item(199, 4).Value = 44488
$ws.Cells.Item(199, 10).Value = 1200
$ws.Cells.Item(199, 11).Value = 1000
$ws.Cells.Item(199, 12).Value = 1200
$ws.Cells.Item(199, 13).Value = 1100
$ws.Cells.Item(199, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(199, 15).Value = 'Región del Maule'
$ws.Cells.Item(199, 16).Value = 220
$ws.Cells.Item(199, 17).Value = 5
$ws.Cells.Item(200, 10).Value = 120
$ws.Cells.Item(200, 11).Value = 8000
$ws.Cells.Item(200, 12).Value = 9000
$ws.Cells.Item(200, 13).Value = 8500
$ws.Cells.Item(200, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(200, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(200, 16).Value = 567
$ws.Cells.Item(200, 17).Value = 15
$ws.Cells.Item(201, 4).Value = 44341
$ws.Cells.Item(201, 10).Value = 1200
$ws.Cells.Item(202, 4).Value = 44452
$ws.Cells.Item(202, 10).Value = 500
$ws.Cells.Item(203, 4).Value = 44473
$ws.Cells.Item(203, 10).Value = 750
$ws.Cells.Item(204, 4).Value = 44298
$ws.Cells.Item(205, 4).Value = 44482
$ws.Cells.Item(205, 10).Value = 250
$ws.Cells.Item(206, 4).Value = 44294
$ws.Cells.Item(206, 10).Value = 750
$ws.Cells.Item(207, 4).Value = 44305
$ws.Cells.Item(207, 10).Value = 250
$ws.Cells.Item(208, 4).Value = 44301
$ws.Cells.Item(208, 10).Value = 500
$ws.Cells.Item(208, 11).Value = 1000
$ws.Cells.Item(208, 13).Value = 1000
$ws.Cells.Item(208, 16).Value = 200
$ws.Cells.Item(209, 4).Value = 44278
$ws.Cells.Item(209, 9).Value = 'Primera'
$ws.Cells.Item(209, 10).Value = 1000
$ws.Cells.Item(209, 11).Value = 850
$ws.Cells.Item(209, 13).Value = 925
$ws.Cells.Item(209, 16).Value = 185
$ws.Cells.Item(210, 4).Value = 44466
$ws.Cells.Item(210, 9).Value = 'Segunda'
$ws.Cells.Item(210, 10).Value = 500
$ws.Cells.Item(210, 11).Value = 1000
$ws.Cells.Item(210, 12).Value = 1000
$ws.Cells.Item(210, 13).Value = 1000
$ws.Cells.Item(210, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(210, 15).Value = 'Región del Maule'
$ws.Cells.Item(210, 16).Value = 200
$ws.Cells.Item(210, 17).Value = 5
$ws.Cells.Item(211, 10).Value = 100
$ws.Cells.Item(211, 11).Value = 10000
$ws.Cells.Item(211, 12).Value = 10000
$ws.Cells.Item(211, 13).Value = 10000
$ws.Cells.Item(211, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(211, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(211, 16).Value = 667
$ws.Cells.Item(211, 17).Value = 15
$ws.Cells.Item(212, 4).Value = 44392
$ws.Cells.Item(212, 10).Value = 1200
$ws.Cells.Item(213, 4).Value = 44322
$ws.Cells.Item(213, 10).Value = 500
$ws.Cells.Item(213, 11).Value = 1000
$ws.Cells.Item(213, 12).Value = 1000
$ws.Cells.Item(213, 13).Value = 1000
$ws.Cells.Item(213, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(213, 15).Value = 'Región del Maule'
$ws.Cells.Item(213, 16).Value = 200
$ws.Cells.Item(213, 17).Value = 5
$ws.Cells.Item(214, 10).Value = 150
$ws.Cells.Item(214, 11).Value = 10000
$ws.Cells.Item(214, 12).Value = 10000
$ws.Cells.Item(214, 13).Value = 10000
$ws.Cells.Item(214, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(214, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(214, 16).Value = 667
$ws.Cells.Item(214, 17).Value = 15
$ws.Cells.Item(215, 4).Value = 44236
$ws.Cells.Item(215, 10).Value = 1200
$ws.Cells.Item(215, 11).Value = 1000
$ws.Cells.Item(215, 12).Value = 1000
$ws.Cells.Item(215, 13).Value = 1000
$ws.Cells.Item(215, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(215, 15).Value = 'Región del Maule'
$ws.Cells.Item(215, 16).Value = 200
$ws.Cells.Item(215, 17).Value = 5
$ws.Cells.Item(216, 10).Value = 120
$ws.Cells.Item(216, 11).Value = 8000
$ws.Cells.Item(216, 12).Value = 8000
$ws.Cells.Item(216, 13).Value = 8000
$ws.Cells.Item(216, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(216, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(216, 16).Value = 533
$ws.Cells.Item(216, 17).Value = 15
$ws.Cells.Item(217, 4).Value = 44162
$ws.Cells.Item(217, 10).Value = 1200
$ws.Cells.Item(218, 4).Value = 44434
$ws.Cells.Item(219, 4).Value = 44532
$ws.Cells.Item(220, 4).Value = 44578
$ws.Cells.Item(220, 11).Value = 1000
$ws.Cells.Item(220, 12).Value = 1000
$ws.Cells.Item(220, 13).Value = 1000
$ws.Cells.Item(220, 16).Value = 200
$ws.Cells.Item(221, 4).Value = 44490
$ws.Cells.Item(221, 11).Value = 1200
$ws.Cells.Item(221, 12).Value = 1200
$ws.Cells.Item(221, 13).Value = 1200
$ws.Cells.Item(221, 16).Value = 240
$ws.Cells.Item(222, 4).Value = 44427
$ws.Cells.Item(222, 10).Value = 500
$ws.Cells.Item(222, 12).Value = 1000
$ws.Cells.Item(222, 13).Value = 1000
$ws.Cells.Item(222, 16).Value = 200
$ws.Cells.Item(223, 4).Value = 44491
$ws.Cells.Item(223, 10).Value = 1200
$ws.Cells.Item(223, 11).Value = 1000
$ws.Cells.Item(223, 12).Value = 1200
$ws.Cells.Item(223, 13).Value = 1100
$ws.Cells.Item(223, 16).Value = 220
$ws.Cells.Item(224, 4).Value = 44266
$ws.Cells.Item(224, 10).Value = 500
$ws.Cells.Item(224, 11).Value = 800
$ws.Cells.Item(224, 13).Value = 900
$ws.Cells.Item(224, 16).Value = 180
$ws.Cells.Item(225, 4).Value = 44533
$ws.Cells.Item(225, 10).Value = 1100
$ws.Cells.Item(226, 4).Value = 44264
$ws.Cells.Item(226, 10).Value = 1000
$ws.Cells.Item(226, 11).Value = 1000
$ws.Cells.Item(226, 12).Value = 1000
$ws.Cells.Item(226, 13).Value = 1000
$ws.Cells.Item(226, 16).Value = 200
$ws.Cells.Item(227, 4).Value = 44494
$ws.Cells.Item(227, 11).Value = 1200
$ws.Cells.Item(227, 12).Value = 1200
$ws.Cells.Item(227, 13).Value = 1200
$ws.Cells.Item(227, 16).Value = 240
$ws.Cells.Item(228, 4).Value = 44571
$ws.Cells.Item(228, 10).Value = 500
$ws.Cells.Item(228, 11).Value = 1000
$ws.Cells.Item(228, 12).Value = 1000
$ws.Cells.Item(228, 13).Value = 1000
$ws.Cells.Item(228, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(228, 15).Value = 'Región del Maule'
$ws.Cells.Item(228, 16).Value = 200
$ws.Cells.Item(228, 17).Value = 5
$ws.Cells.Item(229, 10).Value = 120
$ws.Cells.Item(229, 11).Value = 10000
$ws.Cells.Item(229, 12).Value = 10000
$ws.Cells.Item(229, 13).Value = 10000
$ws.Cells.Item(229, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(229, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(229, 16).Value = 667
$ws.Cells.Item(229, 17).Value = 15
$ws.Cells.Item(230, 4).Value = 44390
$ws.Cells.Item(230, 11).Value = 1000
$ws.Cells.Item(230, 13).Value = 1000
$ws.Cells.Item(230, 16).Value = 200
$ws.Cells.Item(231, 4).Value = 44481
$ws.Cells.Item(231, 10).Value = 1200
$ws.Cells.Item(231, 11).Value = 900
$ws.Cells.Item(231, 12).Value = 1000
$ws.Cells.Item(231, 13).Value = 950
$ws.Cells.Item(231, 16).Value = 190
$ws.Cells.Item(232, 4).Value = 44277
$ws.Cells.Item(232, 11).Value = 850
$ws.Cells.Item(232, 12).Value = 850
$ws.Cells.Item(232, 13).Value = 850
$ws.Cells.Item(232, 15).Value = 'Región del Maule'
$ws.Cells.Item(232, 16).Value = 170
$ws.Cells.Item(233, 4).Value = 44525
$ws.Cells.Item(233, 10).Value = 500
$ws.Cells.Item(233, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(234, 4).Value = 44327
$ws.Cells.Item(234, 10).Value = 1250
$ws.Cells.Item(235, 4).Value = 44354
$ws.Cells.Item(235, 10).Value = 250
$ws.Cells.Item(235, 11).Value = 1000
$ws.Cells.Item(235, 12).Value = 1000
$ws.Cells.Item(235, 13).Value = 1000
$ws.Cells.Item(235, 16).Value = 200
$ws.Cells.Item(236, 4).Value = 44462
$ws.Cells.Item(236, 10).Value = 500
$ws.Cells.Item(236, 11).Value = 1200
$ws.Cells.Item(236, 12).Value = 1200
$ws.Cells.Item(236, 13).Value = 1200
$ws.Cells.Item(236, 16).Value = 240
$ws.Cells.Item(237, 4).Value = 44312
$ws.Cells.Item(237, 10).Value = 250

# --- Add new row 238 (new entry, previously row 237 data shifted into it) ---
$ws.Cells.Item(238, 1).Value = 4
$ws.Cells.Item(238, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(238, 3).Value = 'Los Lagos'
$ws.Cells.Item(238, 4).Value = 44511
$ws.Cells.Item(238, 5).Value = 10
$ws.Cells.Item(238, 6).Value = 100114014
$ws.Cells.Item(238, 7).Value = 'Betarraga'
$ws.Cells.Item(238, 8).Value = 'Sin especificar'
$ws.Cells.Item(238, 9).Value = 'Primera'
$ws.Cells.Item(238, 10).Value = 500
$ws.Cells.Item(238, 11).Value = 1000
$ws.Cells.Item(238, 12).Value = 1000
$ws.Cells.Item(238, 13).Value = 1000
$ws.Cells.Item(238, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(238, 15).Value = 'Región del Maule'
$ws.Cells.Item(238, 16).Value = 200
$ws.Cells.Item(238, 17).Value = 5
$ws.Cells.Item(238, 18).Value = 'Hortaliza'

# Preserve date style/number format on the new row 238 Fecha cell (column D), matching other date cells
$ws.Cells.Item(238, 4).NumberFormat = $ws.Cells.Item(237, 4).NumberFormat()
